# Generate Report for Handoff
# The "2ddde282-..." file was handed off again: its "Latest Handoff Datetime"
# is refreshed on both the zh-cn and the de-de status sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Row 4 corresponds to 2ddde282-c7ce-4e0b-9f58-529427073233 in both sheets.
$zhcn.Range("D4").Value = "2016-03-02 14:14:26"
$dede.Range("D4").Value = "2016-03-02 14:14:36"
